$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 379.44446
$ws.Range("I19").Value = 220.8
$ws.Range("J19").Value = 440.46155
$ws.Range("K19").Value = 220.8
$ws.Range("L19").Value = 440.46155
$ws.Range("M19").Value = -45.80000000000001
$ws.Range("N19").Value = -790.46155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 4678.8
$ws.Range("J34").Value = 14100
$ws.Range("L34").Value = 14100
$ws.Range("N34").Value = -14506

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H36").Value = 4678.8
$ws.Range("J36").Value = 14100
$ws.Range("L36").Value = 14100
$ws.Range("N36").Value = -15530

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 2035.5333
$ws.Range("I127").Value = 648.1429000000001
$ws.Range("J127").Value = 3249.5
$ws.Range("K127").Value = 1944.4287
$ws.Range("L127").Value = 9748.5
$ws.Range("M127").Value = 3015.5713
$ws.Range("N127").Value = -19668.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1055.127
$ws.Range("I129").Value = 463.33334
$ws.Range("J129").Value = 1117.421
$ws.Range("K129").Value = 1390.00002
$ws.Range("L129").Value = 3352.263
$ws.Range("M129").Value = 3609.99998
$ws.Range("N129").Value = -13352.263

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2695.2632
$ws.Range("I132").Value = 2137.375
$ws.Range("J132").Value = 5670.6665
$ws.Range("K132").Value = 6412.125
$ws.Range("L132").Value = 17011.9995
$ws.Range("M132").Value = -3882.125
$ws.Range("N132").Value = -22071.9995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2799.45
$ws.Range("I141").Value = 2160.7693
$ws.Range("J141").Value = 3985.5715
$ws.Range("K141").Value = 6482.3079
$ws.Range("L141").Value = 11956.7145
$ws.Range("M141").Value = -1302.3079
$ws.Range("N141").Value = -22316.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8005.609
$ws.Range("I32").Value = 5572.5083
$ws.Range("K32").Value = 5572.5083
$ws.Range("M32").Value = -5285.5083

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 53579.5
$ws.Range("J52").Value = 53579.5
$ws.Range("L52").Value = 53579.5
$ws.Range("N52").Value = -54215.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1575.35
$ws.Range("I74").Value = 1203.92
$ws.Range("J74").Value = 2194.4
$ws.Range("K74").Value = 1203.92
$ws.Range("L74").Value = 2194.4
$ws.Range("M74").Value = -329.9200000000001
$ws.Range("N74").Value = -3942.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1575.35
$ws.Range("I77").Value = 1203.92
$ws.Range("J77").Value = 2194.4
$ws.Range("K77").Value = 6019.6
$ws.Range("L77").Value = 10972
$ws.Range("M77").Value = -1651.6
$ws.Range("N77").Value = -19708

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11656.917
$ws.Range("I20").Value = 1141.2667
$ws.Range("J20").Value = 29183
$ws.Range("K20").Value = 1141.2667
$ws.Range("L20").Value = 29183
$ws.Range("M20").Value = -894.2666999999999
$ws.Range("N20").Value = -29677

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1097.591
$ws.Range("I99").Value = 782.6429000000001
$ws.Range("J99").Value = 1648.75
$ws.Range("K99").Value = 782.6429000000001
$ws.Range("L99").Value = 1648.75
$ws.Range("M99").Value = 715.3570999999999
$ws.Range("N99").Value = -4644.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 907.94116
$ws.Range("I107").Value = 756.1923
$ws.Range("J107").Value = 1401.125
$ws.Range("K107").Value = 756.1923
$ws.Range("L107").Value = 1401.125
$ws.Range("M107").Value = 1163.8077
$ws.Range("N107").Value = -5241.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3431
$ws.Range("I31").Value = 2422.9524
$ws.Range("J31").Value = 4113.871
$ws.Range("K31").Value = 2422.9524
$ws.Range("L31").Value = 4113.871
$ws.Range("M31").Value = -2127.9524
$ws.Range("N31").Value = -4703.871

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3431
$ws.Range("I34").Value = 2422.9524
$ws.Range("J34").Value = 4113.871
$ws.Range("K34").Value = 2422.9524
$ws.Range("L34").Value = 4113.871
$ws.Range("M34").Value = -2220.9524
$ws.Range("N34").Value = -4517.871

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1585.3636
$ws.Range("I105").Value = 1726.6
$ws.Range("J105").Value = 1282.7142
$ws.Range("K105").Value = 1726.6
$ws.Range("L105").Value = 1282.7142
$ws.Range("M105").Value = 20.40000000000009
$ws.Range("N105").Value = -4776.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 585.0417
$ws.Range("I107").Value = 226.125
$ws.Range("J107").Value = 1302.875
$ws.Range("K107").Value = 226.125
$ws.Range("L107").Value = 1302.875
$ws.Range("M107").Value = 1693.875
$ws.Range("N107").Value = -5142.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3409.1724
$ws.Range("I134").Value = 3578.68
$ws.Range("J134").Value = 2349.75
$ws.Range("K134").Value = 10736.04
$ws.Range("L134").Value = 7049.25
$ws.Range("M134").Value = -8201.039999999999
$ws.Range("N134").Value = -12119.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4591013
$ws.Range("J12").Value = 71524.42999999999
$ws.Range("L12").Value = 214573.29
$ws.Range("N12").Value = -214919.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1099.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1099.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3298.5
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3636.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 489.30435
$ws.Range("I40").Value = 197.57895
$ws.Range("J40").Value = 1875
$ws.Range("K40").Value = 790.3158
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -721.3158
$ws.Range("N40").Value = -7638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 713.5
$ws.Range("J107").Value = 785.4545000000001
$ws.Range("L107").Value = 2356.3635
$ws.Range("N107").Value = -6196.3635

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1429099.8
$ws.Range("I113").Value = 2000529
$ws.Range("J113").Value = 588762.7
$ws.Range("K113").Value = 6001587
$ws.Range("L113").Value = 1766288.1
$ws.Range("M113").Value = -5999417
$ws.Range("N113").Value = -1770628.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2066.1428
$ws.Range("I132").Value = 1658.1111
$ws.Range("J132").Value = 2800.6
$ws.Range("K132").Value = 14922.9999
$ws.Range("L132").Value = 25205.4
$ws.Range("M132").Value = -12392.9999
$ws.Range("N132").Value = -30265.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 377936.4
$ws.Range("I102").Value = 628677.7
$ws.Range("J102").Value = 1824.4445
$ws.Range("K102").Value = 628677.7
$ws.Range("L102").Value = 1824.4445
$ws.Range("M102").Value = -627055.7
$ws.Range("N102").Value = -5068.4445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3934.6226
$ws.Range("I126").Value = 5488.9585
$ws.Range("J126").Value = 2648.276
$ws.Range("K126").Value = 16466.8755
$ws.Range("L126").Value = 7944.828
$ws.Range("M126").Value = -13996.8755
$ws.Range("N126").Value = -12884.828

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 91001864
$ws.Range("I68").Value = 113055.664
$ws.Range("K68").Value = 113055.664
$ws.Range("M68").Value = -112306.664

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 91001864
$ws.Range("I71").Value = 113055.664
$ws.Range("K71").Value = 565278.3200000001
$ws.Range("M71").Value = -561534.3200000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 22400
$ws.Range("J119").Value = 22400
$ws.Range("L119").Value = 22400
$ws.Range("N119").Value = -32076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 968.375
$ws.Range("I96").Value = 1013.5
$ws.Range("J96").Value = 833
$ws.Range("K96").Value = 833
$ws.Range("L96").Value = 833
$ws.Range("M96").Value = 359.5
$ws.Range("N96").Value = -3579
